$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain text in the source workbook, even when the text
# looks like a simple decimal number (e.g. "596.20"). Force those specific cells to stay
# text-formatted before writing the new value so Excel does not silently reinterpret them
# as numeric cells (which would also normalize "596.20" -> 596.2, "1.00" -> 1, etc.).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "67.329.30"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.675.42"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "596.20"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "162.89"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "2.675.03"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "5.20"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "27.76"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "3.164.98"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "67.234.22"
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "2.662.93"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").Value = "11.62"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").Value = "362.42"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("D23").Value = "4.78"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("E24").Value = "  -3.96%  "
$ws.Range("D25").Value = "71.93"
$ws.Range("E25").Value = "  -4.85%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "2.816.11"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "547.66"
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -5.41%  "
$ws.Range("D38").Value = "19.53"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "156.62"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").Value = "17.91"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "40.32"
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("D47").Value = "0.0₆0300"
$ws.Range("E47").Value = "  -5.58%  "
$ws.Range("D48").Value = "0.584"
$ws.Range("E48").Value = "  -4.10%  "
$ws.Range("D49").Value = "152.54"
$ws.Range("E49").Value = "  -4.40%  "
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("E51").Value = "  -2.67%  "
